# Update ERC721F data with latest OpenZeppelin version
# Cell value updates on the three data sheets (columns J:L, rows 4-10).
# Charts on the "Graphs" sheet reference these cells indirectly through
# helper formulas (P/Q/R, Z/AA columns), so once the source values change
# and the workbook recalculates, the dependent formulas pick up the new
# numbers automatically.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: "Optimizer Disabled"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Optimizer Disabled")
$ws1.Range("J4").Value = 74788
$ws1.Range("K4").Value = 307637
$ws1.Range("L4").Value = 2636590
$ws1.Range("J5").Value = 58022
$ws1.Range("K5").Value = 62822
$ws1.Range("J6").Value = 58237
$ws1.Range("K6").Value = 63037
$ws1.Range("K7").Value = 187607
$ws1.Range("L7").Value = 192407
$ws1.Range("K8").Value = 189536
$ws1.Range("L8").Value = 194336
$ws1.Range("L9").Value = 1197308
$ws1.Range("L10").Value = 1206968

# ---------------------------------------------------------------------
# Sheet: "Runs 1000 - Optimizer Enabled"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Runs 1000 - Optimizer Enabled")
$ws3.Range("J4").Value = 74049
$ws3.Range("K4").Value = 301786
$ws3.Range("L4").Value = 2579619
$ws3.Range("J5").Value = 57227
$ws3.Range("K5").Value = 62027
$ws3.Range("J6").Value = 57296
$ws3.Range("K6").Value = 62096
$ws3.Range("K7").Value = 175621
$ws3.Range("L7").Value = 180421
$ws3.Range("K8").Value = 175100
$ws3.Range("L8").Value = 179900
$ws3.Range("L9").Value = 1020982
$ws3.Range("L10").Value = 996392
$ws3.Range("H15").Select() | Out-Null

# Re-select/re-home sheet1 after editing, so its lingering selection
# matches the target (Excel updates the active selection as cells are
# written to, so re-select explicitly).
$ws1.Range("H12").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet: "Runs 200 - Optimizer Enabled" (becomes the active sheet/tab)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Runs 200 - Optimizer Enabled")
$ws2.Range("J4").Value = 74061
$ws2.Range("K4").Value = 301906
$ws2.Range("L4").Value = 2580819
$ws2.Range("J5").Value = 57239
$ws2.Range("K5").Value = 62039
$ws2.Range("J6").Value = 57308
$ws2.Range("K6").Value = 62108
$ws2.Range("K7").Value = 175741
$ws2.Range("L7").Value = 180541
$ws2.Range("K8").Value = 175220
$ws2.Range("L8").Value = 180020
$ws2.Range("L9").Value = 1021582
$ws2.Range("L10").Value = 996992

$ws2.Activate() | Out-Null
$ws2.Range("K10").Select() | Out-Null
